$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) hold text-formatted numbers such as
# "37.886.61" (thousands separated with dots) or trailing-zero decimals
# (e.g. "0.0820", "4.10"). Force the cell format to Text right before
# writing each one so Excel does not reinterpret/round them as numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.886.61"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.047.04"
$ws.Range("E3").Value = "  -0.63%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.18"
$ws.Range("E5").Value = "  -0.57%  "

# Row 6
$ws.Range("E6").Value = "  -1.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.05"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  -2.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0820"
$ws.Range("E10").Value = "  -1.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.17%  "

# Row 12
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.350.69"
$ws.Range("E12").Value = "  -0.57%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.66"
$ws.Range("E13").Value = "  -1.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.06"
$ws.Range("E14").Value = "  -0.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.775"
$ws.Range("E15").Value = "  +1.53%  "

# Row 16
$ws.Range("E16").Value = "  -2.71%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.059.10"
$ws.Range("E17").Value = "  +0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.859.75"
$ws.Range("E18").Value = "  -0.96%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.62"
$ws.Range("E19").Value = "  -0.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  -5.61%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -1.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.96"
$ws.Range("E22").Value = "  -0.72%  "

# Row 23
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +0.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +2.29%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.94"
$ws.Range("E26").Value = "  +1.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.35"
$ws.Range("E27").Value = "  +0.80%  "

# Row 28
$ws.Range("E28").Value = "  -2.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.84"
$ws.Range("E29").Value = "  -0.89%  "

# Row 30
$ws.Range("E30").Value = "  -1.89%  "

# Row 31
$ws.Range("E31").Value = "  -0.17%  "

# Row 32
$ws.Range("E32").Value = "  +8.20%  "

# Row 33
$ws.Range("E33").Value = "  -2.84%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("E34").Value = "  -1.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0604"
$ws.Range("E35").Value = "  -0.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.59"
$ws.Range("E36").Value = "  +3.44%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  +1.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("E38").Value = "  +4.87%  "

# Row 39
$ws.Range("E39").Value = "  +0.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.04"
$ws.Range("E40").Value = "  +6.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.545.19"
$ws.Range("E41").Value = "  +1.13%  "

# Row 42
$ws.Range("E42").Value = "  +0.21%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.34"
$ws.Range("E43").Value = "  -1.78%  "

# Row 44
$ws.Range("E44").Value = "  -1.64%  "

# Row 45
$ws.Range("E45").Value = "  -2.40%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.10"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48
$ws.Range("E48").Value = "  -0.67%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  -0.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.04"
$ws.Range("E50").Value = "  -0.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.238.37"
$ws.Range("E51").Value = "  -0.64%  "
